$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 175949
$ws.Range("C4").Value = 165920
$ws.Range("C5").Value = 10029
$ws.Range("C8").Value = 64.69
